# Inserts one new weekly record (two rows: "Primera" and "Segunda" quality
# grades) for Apio into the historical data table, right before the current
# row 53. This pushes the existing rows 53:128 down to 55:130 and grows the
# sheet's used range from A1:R128 to A1:R130.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows by inserting above the current row 53.
$ws.Rows("53:54").Insert()

# New row 53: Apio, Americana (o), Primera
$ws.Range("A53").Value = 9
$ws.Range("B53").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C53").Value = "Metropolitana"
$ws.Range("D53").Value = 44467
$ws.Range("E53").Value = 13
$ws.Range("F53").Value = 100112017
$ws.Range("G53").Value = "Apio"
$ws.Range("H53").Value = "Americana (o)"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 61
$ws.Range("K53").Value = 7000
$ws.Range("L53").Value = 8000
$ws.Range("M53").Value = 7492
$ws.Range("N53").Value = "`$/docena de matas"
$ws.Range("O53").Value = "Región de Coquimbo"
$ws.Range("P53").Value = 1249
$ws.Range("Q53").Value = 6
$ws.Range("R53").Value = "Hortaliza"

# New row 54: Apio, Americana (o), Segunda
$ws.Range("A54").Value = 9
$ws.Range("B54").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C54").Value = "Metropolitana"
$ws.Range("D54").Value = 44467
$ws.Range("E54").Value = 13
$ws.Range("F54").Value = 100112017
$ws.Range("G54").Value = "Apio"
$ws.Range("H54").Value = "Americana (o)"
$ws.Range("I54").Value = "Segunda"
$ws.Range("J54").Value = 34
$ws.Range("K54").Value = 5000
$ws.Range("L54").Value = 6000
$ws.Range("M54").Value = 5500
$ws.Range("N54").Value = "`$/docena de matas"
$ws.Range("O54").Value = "Región de Coquimbo"
$ws.Range("P54").Value = 917
$ws.Range("Q54").Value = 6
$ws.Range("R54").Value = "Hortaliza"
